$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260, shifting existing rows 260-298 down to 261-299.
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with the new record's data.
$ws.Range("A260").Value2 = 9
$ws.Range("B260").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C260").Value2 = "Metropolitana"
$ws.Range("D260").Value2 = 44776
$ws.Range("D260").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E260").Value2 = 13
$ws.Range("F260").Value2 = 300000001
$ws.Range("G260").Value2 = "Rabanito"
$ws.Range("H260").Value2 = "Sin especificar"
$ws.Range("I260").Value2 = "Primera"
$ws.Range("J260").Value2 = 7000
$ws.Range("K260").Value2 = 2500
$ws.Range("L260").Value2 = 3000
$ws.Range("M260").Value2 = 2750
$ws.Range("N260").Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Range("O260").Value2 = "Provincia de Chacabuco"
$ws.Range("P260").Value2 = 28
$ws.Range("Q260").Value2 = 100
$ws.Range("R260").Value2 = "Hortaliza"
